$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: header / first record (A1 changes from "gv4" to "gv1") ---
$ws.Range("A1").Value = "gv1"
$ws.Range("B1").Value = "giảng viên 1"
$ws.Range("C1").Value = "gv1@gmail.com"
$ws.Range("D1").Value = "cntt"

# --- Fill columns A and B first for the next few rows (matches original
#     authoring order so shared-string indices line up), then C/D ---
$ws.Range("A2").Value = "gv2"
$ws.Range("A3").Value = "gv3"
$ws.Range("B2").Value = "giảng viên 2"
$ws.Range("C2").Value = "gv2@gmail.com"
$ws.Range("D2").Value = "cntt"
$ws.Range("B3").Value = "giảng viên 3"
$ws.Range("C3").Value = "gv1@gmail.com"
$ws.Range("D3").Value = "cntt"

$ws.Range("A4").Value = "gv4"
$ws.Range("B4").Value = "giảng viên 4"
$ws.Range("C4").Value = "gv2@gmail.com"
$ws.Range("D4").Value = "cntt"

$ws.Range("A5").Value = "gv5"
$ws.Range("B5").Value = "giảng viên 5"
$ws.Range("C5").Value = "gv1@gmail.com"
$ws.Range("D5").Value = "cntt"

$ws.Range("A6").Value = "gv6"
$ws.Range("B6").Value = "giảng viên 6"
$ws.Range("C6").Value = "gv2@gmail.com"
$ws.Range("D6").Value = "cntt"

$ws.Range("A7").Value = "gv7"
$ws.Range("B7").Value = "giảng viên 7"
$ws.Range("C7").Value = "gv1@gmail.com"
$ws.Range("D7").Value = "cntt"

$ws.Range("A8").Value = "gv8"
$ws.Range("B8").Value = "giảng viên 8"
$ws.Range("C8").Value = "gv2@gmail.com"
$ws.Range("D8").Value = "cntt"

$ws.Range("A9").Value = "gv9"
$ws.Range("B9").Value = "giảng viên 9"
$ws.Range("C9").Value = "gv1@gmail.com"
$ws.Range("D9").Value = "cntt"

$ws.Range("A10").Value = "gv10"
$ws.Range("B10").Value = "giảng viên 10"
$ws.Range("C10").Value = "gv2@gmail.com"
$ws.Range("D10").Value = "cntt"

$ws.Range("A11").Value = "gv11"
$ws.Range("B11").Value = "giảng viên 11"
$ws.Range("C11").Value = "gv1@gmail.com"
$ws.Range("D11").Value = "cntt"

$ws.Range("A12").Value = "gv12"
$ws.Range("B12").Value = "giảng viên 12"
$ws.Range("C12").Value = "gv2@gmail.com"
$ws.Range("D12").Value = "cntt"

$ws.Range("A13").Value = "gv13"
$ws.Range("B13").Value = "giảng viên 13"
$ws.Range("C13").Value = "gv1@gmail.com"
$ws.Range("D13").Value = "cntt"

$ws.Range("A14").Value = "gv14"
$ws.Range("B14").Value = "giảng viên 14"
$ws.Range("C14").Value = "gv2@gmail.com"
$ws.Range("D14").Value = "cntt"

$ws.Range("A15").Value = "gv15"
$ws.Range("B15").Value = "giảng viên 15"
$ws.Range("C15").Value = "gv1@gmail.com"
$ws.Range("D15").Value = "cntt"

$ws.Range("A16").Value = "gv16"
$ws.Range("B16").Value = "giảng viên 16"
$ws.Range("C16").Value = "gv2@gmail.com"
$ws.Range("D16").Value = "cntt"

$ws.Range("A17").Value = "gv17"
$ws.Range("B17").Value = "giảng viên 17"
$ws.Range("C17").Value = "gv1@gmail.com"
$ws.Range("D17").Value = "cntt"

$ws.Range("A18").Value = "gv18"
$ws.Range("B18").Value = "giảng viên 18"
$ws.Range("C18").Value = "gv2@gmail.com"
$ws.Range("D18").Value = "cntt"

$ws.Range("A19").Value = "gv19"
$ws.Range("B19").Value = "giảng viên 19"
$ws.Range("C19").Value = "gv1@gmail.com"
$ws.Range("D19").Value = "cntt"

# --- Hyperlinks on column C, added in the same order rIds were issued:
#     odd rows (gv1@gmail.com) first, then even rows (gv2@gmail.com) ---
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:gv2@gmail.com")
$ws.Range("C2").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:gv1@gmail.com")
$ws.Range("C3").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("C5"), "mailto:gv1@gmail.com")
$ws.Range("C5").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("C7"), "mailto:gv1@gmail.com")
$ws.Range("C7").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("C9"), "mailto:gv1@gmail.com")
$ws.Range("C9").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("C11"), "mailto:gv1@gmail.com")
$ws.Range("C11").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("C13"), "mailto:gv1@gmail.com")
$ws.Range("C13").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("C15"), "mailto:gv1@gmail.com")
$ws.Range("C15").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("C17"), "mailto:gv1@gmail.com")
$ws.Range("C17").Style = "Hyperlink"

# C19 keeps "gv1@gmail.com" as its visible text but the hyperlink's stored
# display text is "gv2@gmail.com" (mismatch present in the source data).
$ws.Hyperlinks.Add($ws.Range("C19"), "mailto:gv1@gmail.com", [Type]::Missing, [Type]::Missing, "gv2@gmail.com")
$ws.Range("C19").Value = "gv1@gmail.com"
$ws.Range("C19").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:gv2@gmail.com")
$ws.Range("C4").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("C6"), "mailto:gv2@gmail.com")
$ws.Range("C6").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("C8"), "mailto:gv2@gmail.com")
$ws.Range("C8").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("C10"), "mailto:gv2@gmail.com")
$ws.Range("C10").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("C12"), "mailto:gv2@gmail.com")
$ws.Range("C12").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("C14"), "mailto:gv2@gmail.com")
$ws.Range("C14").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("C16"), "mailto:gv2@gmail.com")
$ws.Range("C16").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("C18"), "mailto:gv2@gmail.com")
$ws.Range("C18").Style = "Hyperlink"

# --- View state: select D1:D19 like the source workbook ---
$ws.Activate()
$ws.Range("D1:D19").Select()
